# Apply crypto price/volume updates to Sheet1 (D = Price, E = Volume(1h)).
# Values are stored as literal text in the source workbook; D-column entries
# that look like plain numbers need the cell pre-formatted as Text ("@") before
# the assignment, otherwise Excel's COM layer infers a Number type and the
# OOXML would serialize a <v> numeric cell instead of a shared/inline string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.446.76"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.573.09"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.78"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3727"
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.93"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3392"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07573"
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.145"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.28"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.016"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.959"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "1.575.07"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001122"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.98"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06752"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.300"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.16"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("D24").Value = "22.452.79"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.338"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.693"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.72"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.021"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.68"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "1.751.25"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.062"
$ws.Range("E32").Value = "  +8.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.171"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.986"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.849"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08367"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02493"
$ws.Range("E37").Value = "  -2.28%  "
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2306"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06518"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.466"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.31"
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6219"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.99"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5806"
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.02"
$ws.Range("E48").Value = "  +3.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.066"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("E50").Value = "  -5.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07325"
$ws.Range("E51").Value = "  -0.05%  "
